$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 corresponds to the "Recuperar contraseña" task in the
# Requisitos_de_Proyecto table. Update its status to completed and
# refresh the related notes / result columns (punto 4 doc + excel).
$ws.Range("D13").Value = "Completada"
$ws.Range("G13").Value = 1
$ws.Range("I13").Value = "todo OK, verificación de email, nickname y código"
$ws.Range("J13").Value = "OK"

# Restore the view state (scroll position / active selection) recorded
# in the saved workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K13").Select()
